$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45972, 33.145, 0),
    @(3, 45972.01041666666, 4.643, 3.122),
    @(4, 45972.02083333334, 0.6870000000000001, 1.644),
    @(5, 45972.03125, 3.578, 1.95),
    @(6, 45972.04166666666, 0.003, 25.085),
    @(7, 45972.05208333334, 0.018, 14.168),
    @(8, 45972.0625, 1.85, 2.872),
    @(9, 45972.07291666666, 0, 19.343),
    @(10, 45972.08333333334, 0, 26.124),
    @(11, 45972.09375, 0.152, 6.697),
    @(12, 45972.10416666666, 0, 1.072),
    @(13, 45972.11458333334, 0, 2.029),
    @(14, 45972.125, 1.33, 0.752),
    @(15, 45972.13541666666, 0.097, 1.57),
    @(16, 45972.14583333334, 2.18, 0.5620000000000001),
    @(17, 45972.15625, 0.8110000000000001, 5.381),
    @(18, 45972.16666666666, 0.014, 11.307),
    @(19, 45972.17708333334, 0.001, 10.293),
    @(20, 45972.1875, 0.778, 5.48),
    @(21, 45972.19791666666, 0.061, 10.864),
    @(22, 45972.20833333334, 0.211, 8.537000000000001),
    @(23, 45972.21875, 0.046, 18.425),
    @(24, 45972.22916666666, 0, 36.517),
    @(25, 45972.23958333334, 0.392, 11.859),
    @(26, 45972.25, 0, 33.74),
    @(27, 45972.26041666666, 0, 48.264),
    @(28, 45972.27083333334, 0, 22.214),
    @(29, 45972.28125, 0, 44.75),
    @(30, 45972.29166666666, 0, 58.128),
    @(31, 45972.30208333334, 0, 50.754),
    @(32, 45972.3125, 0, 48.902),
    @(33, 45972.32291666666, 0, 12.569),
    @(34, 45972.33333333334, 0, 23.936),
    @(35, 45972.34375, 0, 22.099),
    @(36, 45972.35416666666, 0.019, 14.428),
    @(37, 45972.36458333334, 0, 30.072),
    @(38, 45972.375, 0.151, 24.135),
    @(39, 45972.38541666666, 0, 15.917),
    @(40, 45972.39583333334, 0, 13.884),
    @(41, 45972.40625, 0, 62.638),
    @(42, 45972.41666666666, 0, 67.321),
    @(43, 45972.42708333334, 0, 14.319),
    @(44, 45972.4375, 0.374, 15.429),
    @(45, 45972.44791666666, 0.393, 5.418),
    @(46, 45972.45833333334, 0.023, 30.065),
    @(47, 45972.46875, 0, 28.932),
    @(48, 45972.47916666666, 2.376, 5.12),
    @(49, 45972.48958333334, 1.675, 2.541),
    @(50, 45972.5, 0.431, 5.108),
    @(51, 45972.51041666666, 2.093, 0.206),
    @(52, 45972.52083333334, 8.56, 0),
    @(53, 45972.53125, 9.212, 0),
    @(54, 45972.54166666666, 0.616, 7.344),
    @(55, 45972.55208333334, 1.811, 2.017),
    @(56, 45972.5625, 1.094, 0.9350000000000001),
    @(57, 45972.57291666666, 0.95, 2.678),
    @(58, 45972.58333333334, 0, 41.501),
    @(59, 45972.59375, 0, 21.909),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
